$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.2281753114861543
$ws.Range("B2").Value = -0.2988997082525781
$ws.Range("A3").Value = -0.4326965787410214
$ws.Range("B3").Value = -0.3423569198510998
$ws.Range("A4").Value = -0.4806731804802525
$ws.Range("B4").Value = -0.4223015327360145
$ws.Range("A5").Value = -0.1697375720733556
$ws.Range("B5").Value = -0.2178148114610044
$ws.Range("A6").Value = -0.1437041167107272
$ws.Range("B6").Value = -0.1270715290878877
$ws.Range("A7").Value = -0.3087157909330267
$ws.Range("B7").Value = -0.2978266703906248
$ws.Range("A8").Value = -0.4588758238577187
$ws.Range("B8").Value = -0.3865040774120232
$ws.Range("A9").Value = -0.575666506424015
$ws.Range("B9").Value = -0.4428253631487147
$ws.Range("A10").Value = -0.4515125164023661
$ws.Range("B10").Value = -0.2602305688054568
$ws.Range("A11").Value = -0.2941110275999199
$ws.Range("B11").Value = -0.291579227651739
$ws.Range("A12").Value = -0.2246411833686304
$ws.Range("B12").Value = -0.1890673678760373
$ws.Range("A13").Value = -0.06854032025181718
$ws.Range("B13").Value = -0.07819663191348619
$ws.Range("A14").Value = -0.2539289094481106
$ws.Range("B14").Value = -0.184774255186958
$ws.Range("A15").Value = -0.1166134280791809
$ws.Range("B15").Value = -0.02791616864122348
$ws.Range("A16").Value = -0.170070281675273
$ws.Range("B16").Value = -0.05626159619898492
$ws.Range("A17").Value = 0.08578589334121202
$ws.Range("B17").Value = 0.1463924458497362
$ws.Range("A18").Value = 0.01934695546148901
$ws.Range("B18").Value = 0.04556771774731427
$ws.Range("A19").Value = 0.02976926871651094
$ws.Range("B19").Value = 0.08822422933785865
$ws.Range("A20").Value = -0.1192795403559119
$ws.Range("B20").Value = -0.0550419988408351
$ws.Range("A21").Value = 0.03906960831486925
$ws.Range("B21").Value = 0.05901722641829715
$ws.Range("A22").Value = 0.0627896640118632
$ws.Range("B22").Value = 0.1241279045274134
$ws.Range("A23").Value = 0.0338726269401972
$ws.Range("B23").Value = 0.03884071518017512
$ws.Range("A24").Value = 0.8113989313042456
$ws.Range("B24").Value = 0.5938562671908888
$ws.Range("A25").Value = 0.1302054731247136
$ws.Range("B25").Value = 0.09281944634686716
$ws.Range("A26").Value = 0.1515642609823832
$ws.Range("B26").Value = 0.1251295574012061
$ws.Range("A27").Value = 0.1055232266952809
$ws.Range("B27").Value = 0.1089181200362979
$ws.Range("A28").Value = 0.2982855054839414
$ws.Range("B28").Value = 0.1882753136055874
$ws.Range("A29").Value = 0.6315833236605232
$ws.Range("B29").Value = 0.5127470668831878
$ws.Range("A30").Value = 0.2014041362207585
$ws.Range("B30").Value = 0.1598295481381244
$ws.Range("A31").Value = 0.02481128159909759
$ws.Range("B31").Value = 0.03514478001957373
$ws.Range("A32").Value = 0.1737050911834538
$ws.Range("B32").Value = 0.1714746687135577
$ws.Range("A33").Value = 0.1182487954320816
$ws.Range("B33").Value = 0.1259171087411987
$ws.Range("A34").Value = 0.08883989179839069
$ws.Range("B34").Value = 0.04951101369534622
$ws.Range("A35").Value = 0.4018459552159545
$ws.Range("B35").Value = 0.2577559197610668
$ws.Range("A36").Value = 0.260730527606143
$ws.Range("B36").Value = 0.1247549647763282
$ws.Range("A37").Value = 0.06585361720161836
$ws.Range("B37").Value = -0.01132238630589366
$ws.Range("A38").Value = 0.3322606920892029
$ws.Range("B38").Value = 0.273073446524114
$ws.Range("A39").Value = -0.07060256033685686
$ws.Range("B39").Value = -0.1582765077064184
$ws.Range("A40").Value = 0.1787789290427507
$ws.Range("B40").Value = 0.1724357899974115
$ws.Range("A41").Value = -0.1169184665546673
$ws.Range("B41").Value = -0.173198470314575
$ws.Range("A42").Value = 0.2418133223432772
$ws.Range("B42").Value = 0.2173282928449247
$ws.Range("A43").Value = 0.1916137478952961
$ws.Range("B43").Value = 0.1530728756710596
$ws.Range("A44").Value = -0.1213111496745746
$ws.Range("B44").Value = -0.09565828784753705
$ws.Range("A45").Value = -0.105815105790153
$ws.Range("B45").Value = -0.1032814672089365
$ws.Range("A46").Value = -0.1894810015761276
$ws.Range("B46").Value = -0.1804892429802131
$ws.Range("A47").Value = -0.1889844693614678
$ws.Range("B47").Value = -0.181158037485951
$ws.Range("A48").Value = -0.2285643802945282
$ws.Range("B48").Value = -0.207368933605906
$ws.Range("A49").Value = -0.2211384405771495
$ws.Range("B49").Value = -0.205019093253885
$ws.Range("A50").Value = -0.1542242326987149
$ws.Range("B50").Value = -0.1444851799490377
$ws.Range("A51").Value = -0.2132275859158974
$ws.Range("B51").Value = -0.2256221463837091
$ws.Range("A52").Value = -0.2132275859158974
$ws.Range("B52").Value = -0.2256221463837091
$ws.Range("A53").Value = -0.2071471597225397
$ws.Range("B53").Value = -0.191734616525058
$ws.Range("A54").Value = -0.207856041027209
$ws.Range("B54").Value = -0.2038544479223625
$ws.Range("A55").Value = -0.1710765624915974
$ws.Range("B55").Value = -0.1552346643470963
$ws.Range("A56").Value = -0.150051687702363
$ws.Range("B56").Value = -0.1457404176251496
$ws.Range("A57").Value = -0.1976795480467335
$ws.Range("B57").Value = -0.1578538221516333
$ws.Range("A58").Value = -0.18036077218508
$ws.Range("B58").Value = -0.2109048882063121
$ws.Range("A59").Value = -0.221485526579741
$ws.Range("B59").Value = -0.2353069304415698
$ws.Range("A60").Value = -0.2094013024649638
$ws.Range("B60").Value = -0.2352997651909498
$ws.Range("A61").Value = -0.2655496591401728
$ws.Range("B61").Value = -0.2272365813879831
$ws.Range("A62").Value = -0.1642320694870018
$ws.Range("B62").Value = -0.1064504400080818
$ws.Range("A63").Value = -0.3436512458246397
$ws.Range("B63").Value = -0.3869883425053283
$ws.Range("A64").Value = -0.2723036396626605
$ws.Range("B64").Value = -0.2638794487875365
$ws.Range("A65").Value = -0.1936561084382792
$ws.Range("B65").Value = -0.2275422745700788
$ws.Range("A66").Value = -0.1196590449759944
$ws.Range("B66").Value = -0.09956497374286856
$ws.Range("A67").Value = -0.0589316749074649
$ws.Range("B67").Value = -0.0601353920044297
